# Create Report in StartBowse.java
#
# Adds a new "Question 8" block (4 answer rows) to the end of the
# "Exam Sample A" sheet, and appends two extra answer options (c/d) to the
# existing "Question 38" block at the end of the "Exam Sample B" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Exam Sample A" (sheet2.xml): append Question 8 after row 52.
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Exam Sample A")

$question8Text = "Question #8`nMatch the following test work products (1-4) with the right description (A-D).`n1. Test suite.`n2. Test case.`n3. Test script.`n4. Test charter.`nA. A set of test scripts to be executed in a specific test run.`nB. A set of instructions for the execution of a test.`nC. Contains expected results.`nD. Documentation of test activities in session-based exploratory testing."

# Blank separator row (53) - touch the cells (without changing formatting)
# so the row is materialized even though every cell in it is empty.
$wsA.Cells.Item(53, 1).Font.Bold = $false
$wsA.Cells.Item(53, 2).Font.Bold = $false
$wsA.Cells.Item(53, 3).Font.Bold = $false

# Header row (54): A is blank, B/C repeat the "Questions"/"Answer" labels.
$wsA.Cells.Item(54, 1).Font.Bold = $false
$wsA.Cells.Item(54, 2).Value = "Questions"
$wsA.Cells.Item(54, 3).Value = "Answer"

# Question 8 + its four answer options (rows 55-59).
$wsA.Cells.Item(55, 1).Value = "Question 8:"
$wsA.Cells.Item(55, 2).Value = $question8Text
$wsA.Cells.Item(55, 3).Value = "Match the following test work products (1-4) with the right description (A-D)."

$wsA.Cells.Item(56, 1).Value = "Question 8:"
$wsA.Cells.Item(56, 2).Value = $question8Text
$wsA.Cells.Item(56, 3).Value = "a) 1A, 2C, 3B, 4D."

$wsA.Cells.Item(57, 1).Value = "Question 8:"
$wsA.Cells.Item(57, 2).Value = $question8Text
$wsA.Cells.Item(57, 3).Value = "b) 1D, 2B, 3A, 4C."

$wsA.Cells.Item(58, 1).Value = "Question 8:"
$wsA.Cells.Item(58, 2).Value = $question8Text
$wsA.Cells.Item(58, 3).Value = "c) 1A, 2C, 3D, 4B."

$wsA.Cells.Item(59, 1).Value = "Question 8:"
$wsA.Cells.Item(59, 2).Value = $question8Text
$wsA.Cells.Item(59, 3).Value = "d) 1D, 2C, 3B, 4A."

# Re-run autofit on the rows holding the multi-line question text so the
# row height is recalculated as "auto" rather than left flagged custom.
$wsA.Range("A55:A59").EntireRow.AutoFit()

# ---------------------------------------------------------------------
# Sheet "Exam Sample B" (sheet3.xml): append two more Question 38
# answer options (c, d) after the existing block ending at row 15.
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Exam Sample B")

$question38Text = "Question #38`nYou are performing system testing of a train reservation system. Based on the test cases performed, you have noticed that the system occasionally reports that no trains are available, although this should actually be the case. You have provided the developers with a summary of the defect and the version of the tested system. They recognize the urgency of the defect and are now waiting for you to provide further details.`nIn addition to the information already provided, the following additional information is given:`n1. Degree of impact (severity) of the defect.`n2. Identification of the test item.`n3. Details of the test environment.`n4. Urgency/priority to fix.`n5. Actual results.`n6. Reference to test case specification.`nWhich of this information is most useful to include in the defect report?"

# Blank separator row (16).
$wsB.Cells.Item(16, 1).Font.Bold = $false
$wsB.Cells.Item(16, 2).Font.Bold = $false
$wsB.Cells.Item(16, 3).Font.Bold = $false

# Header row (17): A is blank, B/C repeat the "Questions"/"Answer" labels.
$wsB.Cells.Item(17, 1).Font.Bold = $false
$wsB.Cells.Item(17, 2).Value = "Questions"
$wsB.Cells.Item(17, 3).Value = "Answer"

# Question 38 restated + its (now six) answer options (rows 18-22).
$wsB.Cells.Item(18, 1).Value = "Question 38:"
$wsB.Cells.Item(18, 2).Value = $question38Text
$wsB.Cells.Item(18, 3).Value = "1. Degree of impact (severity) of the defect."

$wsB.Cells.Item(19, 1).Value = "Question 38:"
$wsB.Cells.Item(19, 2).Value = $question38Text
$wsB.Cells.Item(19, 3).Value = "a) 1, 2, 6"

$wsB.Cells.Item(20, 1).Value = "Question 38:"
$wsB.Cells.Item(20, 2).Value = $question38Text
$wsB.Cells.Item(20, 3).Value = "b) 1, 4, 5, 6"

$wsB.Cells.Item(21, 1).Value = "Question 38:"
$wsB.Cells.Item(21, 2).Value = $question38Text
$wsB.Cells.Item(21, 3).Value = "c) 2, 3, 4, 5"

$wsB.Cells.Item(22, 1).Value = "Question 38:"
$wsB.Cells.Item(22, 2).Value = $question38Text
$wsB.Cells.Item(22, 3).Value = "d) 3, 5, 6"

# Re-run autofit on the rows holding the multi-line question text so the
# row height is recalculated as "auto" rather than left flagged custom.
$wsB.Range("A18:A22").EntireRow.AutoFit()
